$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 14 (currently K_QUALI), shifting rows 14:23 down to 15:24.
$ws.Range("A14:C14").Insert()

# Copy the style/format of the row below (the now-shifted K_QUALI row, row 15) onto
# the newly inserted blank row 14, so it matches the rest of the data rows.
$ws.Range("A15:C15").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the values for the new row 14: K_PRAEV
$ws.Range("A14").Value = "K_PRAEV"
$ws.Range("B14").Value = "Art der Prävention"
$ws.Range("C14").Value = "XXXArt der Prävention"
